$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New team label texts for B3:B27 (B2 becomes blank)
$values = @(
    "    הקבוצה של: שחקן נא, שחקן כ, שחקן מו",
    "    הקבוצה של: שחקן ח, שחקן יב",
    "    הקבוצה של: שחקן יח, שחקן לד",
    "    הקבוצה של: שחקן מ, שחקן נב",
    "    הקבוצה של: שחקן ה, שחקן לו",
    "    הקבוצה של: שחקן לב, שחקן טו",
    "    הקבוצה של: שחקן מה, שחקן כא",
    "    הקבוצה של: שחקן כו, שחקן יג",
    "    הקבוצה של: שחקן מא, שחקן א",
    "    הקבוצה של: שחקן ב, שחקן מג",
    "    הקבוצה של: שחקן כד, שחקן לג",
    "    הקבוצה של: שחקן כח, שחקן י",
    "    הקבוצה של: Sahkan, שחקן ו",
    "    הקבוצה של: שחקן טז, שחקן כה",
    "    הקבוצה של: שחקן לה, שחקן מט",
    "    הקבוצה של: שחקן לא, שחקן מב",
    "    הקבוצה של: שחקן מח, שחקן מד",
    "    הקבוצה של: שחקן ל, שחקן כט",
    "    הקבוצה של: שחקן מז, שחקן יט",
    "    הקבוצה של: שחקן כג, שחקן כז",
    "    הקבוצה של: שחקן יא, שחקן נ",
    "    הקבוצה של: שחקן לט, שחקן כב",
    "    הקבוצה של: שחקן לח, שחקן ז",
    "    הקבוצה של: שחקן יד, שחקן יז",
    "    הקבוצה של: שחקן ט, שחקן לז"
)

# B2 is cleared out (no more text / value)
$ws.Range("B2").Value = $null

# B3..B27 get the new values, in row order
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Row 6-10's label cells pick up the "right aligned" team-label formatting
# (matching B2-B5 / B11-B14), copied from B3 which already carries it.
$ws.Range("B3").Copy()
$ws.Range("B6:B10").PasteSpecial(-4122)  # xlPasteFormats

# Row 14's label cell switches to the plain (non-right-aligned) formatting
# used by the rows below it, copied from B15.
$ws.Range("B15").Copy()
$ws.Range("B14").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# Reset the active selection to B3 (was D9)
$ws.Range("B3").Select()
